$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A31").Value = 25
$ws.Range("B31").Value = "3：17-5：45"
$ws.Range("C31").Value = "第七章结束，明天开始Unix"

$ws.Range("C31").Select()
